$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K" strikeouts column),
# regenerated from Strike# to K per commit message.
$updates = @{
    2  = 0
    4  = 1
    5  = 2
    7  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
